$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "26.024.98"
$ws.Cells.Item(2, 5).Value = "  +0.39%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.643.08"
$ws.Cells.Item(3, 5).Value = "  +0.59%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.55%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "216.42"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.77%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.507"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.90%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.53%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.27%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +1.17%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "19.63"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -0.09%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0795"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.51%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "1.700.42"
$ws.Cells.Item(12, 5).Value = "  +3.78%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "1.872.20"
$ws.Cells.Item(13, 5).Value = "  +0.68%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  +1.26%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.544"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.02%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "0.0₃0766"
$ws.Cells.Item(16, 5).Value = "  +1.24%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +0.74%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "26.112.54"

# Row 19
$ws.Cells.Item(19, 5).Value = "  +0.60%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "193.23"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.01%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  -0.90%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  -0.38%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -0.53%  "

# Row 24
$ws.Cells.Item(24, 2).Value = "Stellar"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.131"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +4.77%  "

# Row 25
$ws.Cells.Item(25, 2).Value = "Toncoin"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "1.80"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.35%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "BinanceUSD"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "1.01"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.89%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "143.88"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.85%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +0.48%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.53"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.32%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +1.19%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.0496"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -0.65%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.27"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.32%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  -0.82%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -3.41%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +1.25%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.905"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.27%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "1.131.58"
$ws.Cells.Item(37, 5).Value = "  -0.53%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -1.93%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -0.54%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +0.06%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +0.09%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "99.42"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.30%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -0.73%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "1.782.12"
$ws.Cells.Item(44, 5).Value = "  +0.75%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  +4.65%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "56.60"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.69%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -0.29%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "7.70"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +1.40%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  +0.18%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "USDD"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.01"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.75%  "
